# Append a new row (row 10) to the receipts sheet with:
#   A10 = "N"
#   B10 = " 11.12"   (leading space, kept as literal text - not a number)
#   C10 = " Not Found"
#
# The workbook's existing rows store numeric-looking values (e.g. B6=" 23.19")
# as plain text, so the new row must match that: assign with a leading
# apostrophe to stop Excel from auto-coercing " 11.12" into the number
# 11.12, then reset the cell style back to Normal so no stray
# quote-prefix/number-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "N"
$ws.Range("B10").Value = "'" + " 11.12"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = " Not Found"
